# IMU fusion decision matrix - restructure the header into a 2-row merged
# block (Criteria/Weight spanning rows 2-3), add a "Total" row, and draw a
# full grid of thin borders around the whole B2:G9 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- 1. Content moves ---------------------------------------------------
# "Criteria" / "Weight" headers move from row 3 up into row 2 (they will be
# vertically merged with row 3 beneath them).
$ws.Range("B2").Value = $ws.Range("B3").Value2
$ws.Range("C2").Value = $ws.Range("C3").Value2
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# New "Total" row label.
$ws.Range("B9").Value = "Total"

# --- 2. Merges ------------------------------------------------------------
$ws.Range("B2:B3").Merge()
$ws.Range("C2:C3").Merge()
$ws.Range("B9:C9").Merge()

# --- 3. Alignment ----------------------------------------------------------
# Center + vertical-center
$ws.Range("B2:C3").HorizontalAlignment = $xlCenter
$ws.Range("B2:C3").VerticalAlignment = $xlCenter
$ws.Range("D3").HorizontalAlignment = $xlCenter
$ws.Range("D3").VerticalAlignment = $xlCenter
$ws.Range("G3").HorizontalAlignment = $xlCenter
$ws.Range("G3").VerticalAlignment = $xlCenter

# Center only
$ws.Range("D2:G2").HorizontalAlignment = $xlCenter
$ws.Range("B9:C9").HorizontalAlignment = $xlCenter

# Center + vertical-center + wrap
$ws.Range("E3:F3").HorizontalAlignment = $xlCenter
$ws.Range("E3:F3").VerticalAlignment = $xlCenter
$ws.Range("E3:F3").WrapText = $true

# Vertical-center only (already present before the edit, kept explicit)
$ws.Range("D5:G5").VerticalAlignment = $xlCenter

# --- 4. Row height for the wrapped sub-header row -------------------------
$ws.Rows.Item(3).RowHeight = 29.15

Write-Host "stage1+2+3 done"
